$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected indices for the dubosson data (rows 2-8, columns B "start" and C "end")
$ws.Range("C2").Value = 58

$ws.Range("B3").Value = 59
$ws.Range("C3").Value = 346

$ws.Range("B4").Value = 347
$ws.Range("C4").Value = 634

$ws.Range("B5").Value = 635
$ws.Range("C5").Value = 923

$ws.Range("B6").Value = 924
$ws.Range("C6").Value = 1211

$ws.Range("B7").Value = 1470
$ws.Range("C7").Value = 1757

$ws.Range("B8").Value = 1758
$ws.Range("C8").Value = 2045

# Update the active selection to match the saved view state
$ws.Range("E10").Select()
